$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.406.88'
$ws.Range("E2").Value = '  -0.83%  '
$ws.Range("D3").Value = '3.515.79'
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''608.83'
$ws.Range("E5").Value = '  -1.58%  '
$ws.Range("D6").Value = '''150.76'
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("D7").Value = '3.514.68'
$ws.Range("E7").Value = '  -1.21%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '''0.481'
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").Value = '''0.139'
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").Value = '''7.03'
$ws.Range("E11").Value = '  +1.00%  '
$ws.Range("D12").Value = '''0.425'
$ws.Range("E12").Value = '  -1.67%  '
$ws.Range("D13").Value = '''0.0000219'
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").Value = '4.109.99'
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").Value = '''31.79'
$ws.Range("E15").Value = '  -0.81%  '
$ws.Range("D16").Value = '3.522.91'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '67.382.88'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '''6.40'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = '''15.27'
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("D21").Value = '''443.08'
$ws.Range("E21").Value = '  -3.29%  '
$ws.Range("D22").Value = '''9.21'
$ws.Range("E22").Value = '  -4.91%  '
$ws.Range("D23").Value = '''0.624'
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("D24").Value = '''77.24'
$ws.Range("E24").Value = '  -0.62%  '
$ws.Range("D25").Value = '''0.0000128'
$ws.Range("E25").Value = '  +10.32%  '
$ws.Range("D26").Value = '3.657.18'
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '''10.13'
$ws.Range("E28").Value = '  -4.90%  '
$ws.Range("D29").Value = '''8.30'
$ws.Range("E29").Value = '  -0.83%  '
$ws.Range("D30").Value = '''2.50'
$ws.Range("E30").Value = '  -2.63%  '
$ws.Range("D31").Value = '''1.55'
$ws.Range("E31").Value = '  -4.97%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +4.21%  '
$ws.Range("D34").Value = '''25.77'
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").Value = '''6.13'
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("D36").Value = '3.506.84'
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").Value = '''1.84'
$ws.Range("E37").Value = '  -3.77%  '
$ws.Range("D38").Value = '''8.03'
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = '''177.18'
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '''1.25'
$ws.Range("E48").Value = '  +3.77%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '''2.59'
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").Value = '''7.58'
$ws.Range("E50").Value = '  -1.81%  '
$ws.Range("D51").Value = '''0.996'
$ws.Range("E51").Value = '  -1.30%  '
